# Refresh the scraped "cryptos" price/volume report (GitHub Actions run).
# Price (column D) and Volume(1h) (column E) are stored as plain text in
# this sheet, not numbers. A handful of the new Price values parse as a
# plain decimal (e.g. "0.4716"), so a bare Range.Value assignment would be
# auto-coerced into a number by Excel; prefixing those with a leading
# apostrophe forces them to stay literal text, matching the original
# formatting. Values that already contain a non-numeric character (a
# second "." thousands separator, a "%" sign, letters, etc.) are
# unambiguous and are assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.194.65'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '1.864.85'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'242.38"
$ws.Range('E5').Value = '  +3.03%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = "'0.4716"
$ws.Range('D8').Value = "'42.67"
$ws.Range('E8').Value = '  -3.29%  '
$ws.Range('D9').Value = "'0.2855"
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('D10').Value = "'0.06462"
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('D11').Value = "'20.71"
$ws.Range('E11').Value = '  -5.05%  '
$ws.Range('D12').Value = "'0.07710"
$ws.Range('E12').Value = '  -3.07%  '
$ws.Range('D13').Value = '1.855.89'
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').Value = "'94.67"
$ws.Range('D15').Value = "'0.6950"
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').Value = "'5.070"
$ws.Range('E16').Value = '  -1.03%  '
$ws.Range('D17').Value = "'267.79"
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '30.176.23'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = "'13.30"
$ws.Range('E19').Value = '  -5.34%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = "'1.001"
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = "'0.000007519"
$ws.Range('E21').Value = '  -2.93%  '
$ws.Range('D22').Value = '2.109.95'
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = "'5.172"
$ws.Range('E24').Value = '  -1.92%  '
$ws.Range('D25').Value = "'6.108"
$ws.Range('E25').Value = '  -1.85%  '
$ws.Range('D26').Value = "'9.297"
$ws.Range('E26').Value = '  -1.35%  '
$ws.Range('D27').Value = "'165.58"
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('D28').Value = "'18.77"
$ws.Range('E28').Value = '  -1.06%  '
$ws.Range('D29').Value = "'1.896"
$ws.Range('E29').Value = '  -3.08%  '
$ws.Range('D30').Value = "'1.379"
$ws.Range('E30').Value = '  +0.56%  '
$ws.Range('D31').Value = "'0.09831"
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('D32').Value = "'1.502"
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('D33').Value = "'4.229"
$ws.Range('E33').Value = '  -2.78%  '
$ws.Range('D34').Value = "'3.999"
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('D35').Value = "'0.04694"
$ws.Range('E35').Value = '  -1.05%  '
$ws.Range('D36').Value = "'1.112"
$ws.Range('E36').Value = '  -2.31%  '
$ws.Range('D37').Value = "'0.6856"
$ws.Range('E37').Value = '  -2.64%  '
$ws.Range('D38').Value = "'2.706"
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').Value = "'0.01837"
$ws.Range('E39').Value = '  -2.14%  '
$ws.Range('D40').Value = "'2.723"
$ws.Range('E40').Value = '  -2.94%  '
$ws.Range('D41').Value = "'6.308"
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('D42').Value = "'70.25"
$ws.Range('E42').Value = '  -3.00%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'0.8416"
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = "'1.000"
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').Value = "'1.884"
$ws.Range('E45').Value = '  -4.03%  '
$ws.Range('D46').Value = "'101.92"
$ws.Range('E46').Value = '  -0.91%  '
$ws.Range('D47').Value = "'0.4055"
$ws.Range('E47').Value = '  -3.16%  '
$ws.Range('D48').Value = "'9.213"
$ws.Range('E48').Value = '  +0.90%  '
$ws.Range('D49').Value = "'7.062"
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('D50').Value = "'932.61"
$ws.Range('E50').Value = '  +1.31%  '
$ws.Range('D51').Value = "'34.54"
$ws.Range('E51').Value = '  -0.22%  '
